$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textAddrs = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D13", "D14", "D15", "D16", "D17", "D19", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D38", "D40", "D41", "D42", "D44", "D45", "D46", "D47", "D49", "D50")
foreach ($addr in $textAddrs) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '29.340.24'
$ws.Range("E2").Value = '  -0.28%  '

$ws.Range("D3").Value = '1.842.08'
$ws.Range("E3").Value = '  -0.67%  '

$ws.Range("D4").Value = '0.9983'
$ws.Range("E4").Value = '  -0.38%  '

$ws.Range("D5").Value = '240.11'
$ws.Range("E5").Value = '  -0.47%  '

$ws.Range("D6").Value = '0.6292'
$ws.Range("E6").Value = '  -0.39%  '

$ws.Range("D7").Value = '0.9993'
$ws.Range("E7").Value = '  -0.14%  '

$ws.Range("D8").Value = '0.07425'
$ws.Range("E8").Value = '  -1.87%  '

$ws.Range("D9").Value = '0.2901'
$ws.Range("E9").Value = '  -0.79%  '

$ws.Range("D10").Value = '24.83'
$ws.Range("E10").Value = '  +1.02%  '

$ws.Range("E11").Value = '  -0.23%  '

$ws.Range("D12").Value = '1.842.45'

$ws.Range("D13").Value = '4.983'
$ws.Range("E13").Value = '  -1.02%  '

$ws.Range("D14").Value = '0.6784'
$ws.Range("E14").Value = '  -0.86%  '

$ws.Range("D15").Value = '0.00001025'
$ws.Range("E15").Value = '  -1.93%  '

$ws.Range("D16").Value = '81.99'
$ws.Range("E16").Value = '  -1.65%  '

$ws.Range("D17").Value = '6.252'
$ws.Range("E17").Value = '  +1.78%  '

$ws.Range("D18").Value = '29.325.31'
$ws.Range("E18").Value = '  -0.32%  '

$ws.Range("D19").Value = '229.10'
$ws.Range("E19").Value = '  -0.44%  '

$ws.Range("E20").Value = '  -0.47%  '

$ws.Range("D21").Value = '0.9996'
$ws.Range("E21").Value = '  -0.15%  '

$ws.Range("D22").Value = '7.436'
$ws.Range("E22").Value = '  -1.00%  '

$ws.Range("D23").Value = '0.9997'
$ws.Range("E23").Value = '  -0.12%  '

$ws.Range("D24").Value = '158.46'
$ws.Range("E24").Value = '  -0.48%  '

$ws.Range("D25").Value = '8.474'
$ws.Range("E25").Value = '  +0.18%  '

$ws.Range("D26").Value = '0.1354'
$ws.Range("E26").Value = '  -3.11%  '

$ws.Range("D27").Value = '17.43'
$ws.Range("E27").Value = '  -1.69%  '

$ws.Range("D28").Value = '0.06545'
$ws.Range("E28").Value = '  +15.05%  '

$ws.Range("D29").Value = '1.448'
$ws.Range("E29").Value = '  +2.46%  '

$ws.Range("D30").Value = '1.488'
$ws.Range("E30").Value = '  +0.43%  '

$ws.Range("D31").Value = '4.068'
$ws.Range("E31").Value = '  -2.03%  '

$ws.Range("D32").Value = '4.064'
$ws.Range("E32").Value = '  +0.14%  '

$ws.Range("D33").Value = '1.838'
$ws.Range("E33").Value = '  +0.61%  '

$ws.Range("E34").Value = '  -1.50%  '

$ws.Range("D35").Value = '0.6962'
$ws.Range("E35").Value = '  -0.40%  '

$ws.Range("D36").Value = '2.572'
$ws.Range("E36").Value = '  -0.66%  '

$ws.Range("E37").Value = '  +1.57%  '

$ws.Range("D38").Value = '2.818'
$ws.Range("E38").Value = '  +1.68%  '

$ws.Range("D39").Value = '1.239.07'
$ws.Range("E39").Value = '  -0.91%  '

$ws.Range("D40").Value = '6.787'
$ws.Range("E40").Value = '  +4.00%  '

$ws.Range("D41").Value = '0.9342'
$ws.Range("E41").Value = '  +2.76%  '

$ws.Range("D42").Value = '0.9991'
$ws.Range("E42").Value = '  -0.14%  '

$ws.Range("D43").Value = '1.985.49'
$ws.Range("E43").Value = '  -1.47%  '

$ws.Range("D44").Value = '100.76'
$ws.Range("E44").Value = '  -0.78%  '

$ws.Range("D45").Value = '65.61'
$ws.Range("E45").Value = '  -0.77%  '

$ws.Range("D46").Value = '0.00000000119'
$ws.Range("E46").Value = '  +2.64%  '

$ws.Range("D47").Value = '7.055'
$ws.Range("E47").Value = '  -1.22%  '

$ws.Range("E48").Value = '  +2.32%  '

$ws.Range("D49").Value = '0.1152'
$ws.Range("E49").Value = '  -1.13%  '

$ws.Range("D50").Value = '8.996'
$ws.Range("E50").Value = '  -0.78%  '

$ws.Range("E51").Value = '  -1.65%  '
